$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Object")

# ---------------------------------------------------------------------------
# "CuTru" object param table (originally rows 5-22): add 3 new param rows
#   - "email"  (string)  -> "Email người đăng ký"
#   - "hoTen"  (string)  -> "Họ tên công dân của Cư Trú"
#   - "diaChi" (string)  -> "Địa chỉ"
# ---------------------------------------------------------------------------

# Step A: insert a row at 13 (just after "ngayHetHan") and fill it as the
# "hoTen" param row.
$ws.Rows.Item(13).Insert()
$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("B9").Copy($ws.Range("B13"))
$ws.Range("C9").Copy($ws.Range("C13"))
$ws.Cells.Item(13, 1).Value2 = "hoTen"
$ws.Cells.Item(13, 2).Value2 = "string"
$ws.Cells.Item(13, 3).Value2 = "Họ tên công dân của Cư Trú"

# Step B: insert another row at 13 (pushes the "hoTen" row down to 14) and
# fill it as the "email" param row.
$ws.Rows.Item(13).Insert()
$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("B9").Copy($ws.Range("B13"))
$ws.Range("C9").Copy($ws.Range("C13"))
$ws.Cells.Item(13, 1).Value2 = "email"
$ws.Cells.Item(13, 2).Value2 = "string"
$ws.Cells.Item(13, 3).Value2 = "Email người đăng ký"

# Step C: insert a row before "daDuyet" (now at row 21 after the two inserts
# above) and fill it as the "diaChi" param row.
$ws.Rows.Item(21).Insert()
$ws.Range("A9").Copy($ws.Range("A21"))
$ws.Range("B9").Copy($ws.Range("B21"))
$ws.Range("C9").Copy($ws.Range("C21"))
$ws.Cells.Item(21, 1).Value2 = "diaChi"
$ws.Cells.Item(21, 2).Value2 = "string"
$ws.Cells.Item(21, 3).Value2 = "Địa chỉ"

# ---------------------------------------------------------------------------
# Defined names referencing the "NguoiDung" / "TrangChu" object headers need
# to follow those blocks, which shifted down by 3 rows (25->28, 46->49).
# ---------------------------------------------------------------------------
$wb.Names.Item("NguoiDungDO").RefersTo = "='Data Object'!`$B`$28"
$wb.Names.Item("TrangChuDO").RefersTo = "='Data Object'!`$B`$49"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the authored edit:
#   - "Quản Lý Cư Trú" is no longer the active tab; its lingering selection
#     moves from F17 to F13.
#   - "Data Object" becomes the active tab, with selection at F14.
# ---------------------------------------------------------------------------
$wsCuTru = $wb.Worksheets.Item("Quản Lý Cư Trú")
$wsCuTru.Activate()
$wsCuTru.Range("F13").Select()

$ws.Activate()
$ws.Range("F14").Select()
